$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Update the date line
$d.Content.Find.Execute("2025-07-08 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-09 Wednesday", 2) | Out-Null

# Row 1
$t.Cell(1, 1).Range.Text = "85÷5=17, 0"
$t.Cell(1, 2).Range.Text = "17÷8=2, 1"
$t.Cell(1, 3).Range.Text = "93÷7=13, 2"
$t.Cell(1, 4).Range.Text = "19÷6=3, 1"
$t.Cell(1, 5).Range.Text = "73÷2=36, 1"

# Row 5
$t.Cell(5, 1).Range.Text = "54÷6=9, 0"
$t.Cell(5, 2).Range.Text = "60÷9=6, 6"
$t.Cell(5, 3).Range.Text = "62÷3=20, 2"
$t.Cell(5, 4).Range.Text = "73÷3=24, 1"
$t.Cell(5, 5).Range.Text = "79÷6=13, 1"

# Row 9
$t.Cell(9, 1).Range.Text = "22÷6=3, 4"
$t.Cell(9, 2).Range.Text = "79÷4=19, 3"
$t.Cell(9, 3).Range.Text = "90÷6=15, 0"
$t.Cell(9, 4).Range.Text = "35÷5=7, 0"
$t.Cell(9, 5).Range.Text = "19÷4=4, 3"

# Row 13
$t.Cell(13, 1).Range.Text = "41÷8=5, 1"
$t.Cell(13, 2).Range.Text = "36÷7=5, 1"
$t.Cell(13, 3).Range.Text = "21÷6=3, 3"
$t.Cell(13, 4).Range.Text = "10÷8=1, 2"
$t.Cell(13, 5).Range.Text = "48÷8=6, 0"

# Row 17
$t.Cell(17, 1).Range.Text = "67÷8=8, 3"
$t.Cell(17, 2).Range.Text = "36÷3=12, 0"
$t.Cell(17, 3).Range.Text = "42÷7=6, 0"
$t.Cell(17, 4).Range.Text = "18÷6=3, 0"
$t.Cell(17, 5).Range.Text = "58÷9=6, 4"
